$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: "Unnamed: 0" -> "ID"
$ws.Range("F1").Value = "ID"

# H2: clear (was "O254")
$ws.Range("H2").Value = ""

# L2: "ACTIVE" -> "INACTIVE"
$ws.Range("L2").Value = "INACTIVE"

# M2: "BAJAJ" -> "IDFC"
$ws.Range("M2").Value = "IDFC"

# N2: "CD" -> "TW"
$ws.Range("N2").Value = "TW"

# H3: clear (was "O281")
$ws.Range("H3").Value = ""

# H5: clear (was "O254")
$ws.Range("H5").Value = ""

# L5: "ACTIVE" -> "INACTIVE"
$ws.Range("L5").Value = "INACTIVE"
